$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 193 updates
$ws.Range("D193").Value = 45075
$ws.Range("N193").Value = 26000
$ws.Range("O193").Value = 27000
$ws.Range("P193").Value = 26625
$ws.Range("S193").Value = 2662

# Row 194 updates
$ws.Range("D194").Value = 45075
$ws.Range("M194").Value = 300
$ws.Range("N194").Value = 23000
$ws.Range("O194").Value = 24000
$ws.Range("P194").Value = 23500
$ws.Range("S194").Value = 2350

# Row 195 updates
$ws.Range("D195").Value = 44741
$ws.Range("K195").Value = 'Hass'
$ws.Range("L195").Value = 'Primera'
$ws.Range("M195").Value = 400
$ws.Range("N195").Value = 15000
$ws.Range("O195").Value = 16000
$ws.Range("P195").Value = 15500
$ws.Range("Q195").Value = '$/bandeja 10 kilos'
$ws.Range("R195").Value = 'Perú'
$ws.Range("S195").Value = 1550
$ws.Range("T195").Value = 10

# Row 196 updates
$ws.Range("D196").Value = 45014
$ws.Range("K196").Value = 'Hass'
$ws.Range("L196").Value = 'Segunda'
$ws.Range("M196").Value = 400
$ws.Range("N196").Value = 30000
$ws.Range("O196").Value = 31000
$ws.Range("P196").Value = 30625
$ws.Range("Q196").Value = '$/bandeja 10 kilos'
$ws.Range("R196").Value = 'Perú'
$ws.Range("S196").Value = 3062
$ws.Range("T196").Value = 10

# Row 197 updates
$ws.Range("D197").Value = 44398
$ws.Range("K197").Value = 'Fuerte'
$ws.Range("M197").Value = 120
$ws.Range("N197").Value = 59000
$ws.Range("O197").Value = 60000
$ws.Range("P197").Value = 59500
$ws.Range("S197").Value = 2380

# Row 198 updates
$ws.Range("D198").Value = 44398
$ws.Range("M198").Value = 100
$ws.Range("N198").Value = 58000
$ws.Range("O198").Value = 59000
$ws.Range("P198").Value = 58500
$ws.Range("S198").Value = 2340

# Row 199 (new row)
$ws.Range("A199").Value = 1
$ws.Range("B199").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C199").Value = 'Arica y Parinacota'
$ws.Range("D199").Value = 44160
$ws.Range("D199").NumberFormat = $ws.Range("D197").NumberFormat
$ws.Range("E199").Value = 15
$ws.Range("F199").Value = 'Fruta'
$ws.Range("G199").Value = 100106
$ws.Range("H199").Value = 'Oleaginosos'
$ws.Range("I199").Value = 100106002
$ws.Range("J199").Value = 'Palta'
$ws.Range("K199").Value = 'Edranol'
$ws.Range("L199").Value = 'Segunda'
$ws.Range("M199").Value = 270
$ws.Range("N199").Value = 68000
$ws.Range("O199").Value = 70000
$ws.Range("P199").Value = 69000
$ws.Range("Q199").Value = '$/caja 25 kilos'
$ws.Range("R199").Value = 'Región de Coquimbo'
$ws.Range("S199").Value = 2760
$ws.Range("T199").Value = 25

# Row 200 (new row)
$ws.Range("A200").Value = 1
$ws.Range("B200").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C200").Value = 'Arica y Parinacota'
$ws.Range("D200").Value = 44356
$ws.Range("D200").NumberFormat = $ws.Range("D197").NumberFormat
$ws.Range("E200").Value = 15
$ws.Range("F200").Value = 'Fruta'
$ws.Range("G200").Value = 100106
$ws.Range("H200").Value = 'Oleaginosos'
$ws.Range("I200").Value = 100106002
$ws.Range("J200").Value = 'Palta'
$ws.Range("K200").Value = 'Fuerte'
$ws.Range("L200").Value = 'Tercera'
$ws.Range("M200").Value = 150
$ws.Range("N200").Value = 64000
$ws.Range("O200").Value = 65000
$ws.Range("P200").Value = 64500
$ws.Range("Q200").Value = '$/caja 25 kilos'
$ws.Range("R200").Value = 'Región de Coquimbo'
$ws.Range("S200").Value = 2580
$ws.Range("T200").Value = 25
